$d = $word.ActiveDocument

# 1. Trim the trailing parenthetical clause from the RSS problem statement.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "100이하의 모든 양의 정수 a와 b중, RSS를 최소화하는 a와 b를 구하여라. (단, a와 b는 모두 100 이하의 양의 정수이다)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "100이하의 모든 양의 정수 a와 b중, RSS를 최소화하는 a와 b를 구하여라.",
    2
)

# 2. Tighten the auto line spacing (256/240 -> 254/240) on every paragraph style.
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    $style = $d.Styles.Item($i)
    $style.ParagraphFormat.LineSpacing = 12.7
}
